$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text so numeric-looking price strings
# (e.g. "1.00", "29.20", "0.520") are not coerced to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.017.47"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "3.307.94"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "588.14"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").Value = "185.57"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("D9").Value = "0.136"
$ws.Range("E9").Value = "  +5.35%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "0.423"
$ws.Range("E11").Value = "  +2.74%  "
$ws.Range("D12").Value = "3.881.68"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "29.20"
$ws.Range("D15").Value = "69.036.51"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("E16").Value = "  +3.79%  "
$ws.Range("D17").Value = "3.314.87"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").Value = "13.70"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("D20").Value = "391.84"
$ws.Range("E20").Value = "  +5.14%  "
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "72.07"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +3.59%  "
$ws.Range("D25").Value = "0.520"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.189"
$ws.Range("E26").Value = "  +4.82%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "9.77"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "5.81"
$ws.Range("E29").Value = "  +3.12%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "2.00"
$ws.Range("E30").Value = "  +1.70%  "
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("D32").Value = "1.32"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("E33").Value = "  +5.29%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("D36").Value = "163.25"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +4.06%  "
$ws.Range("D38").Value = "0.839"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "26.63"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "4.61"
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").Value = "2.63"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").Value = "41.73"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("D44").Value = "0.0694"
$ws.Range("E44").Value = "  +3.46%  "
$ws.Range("D45").Value = "25.49"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "2.639.58"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "343.06"
$ws.Range("E47").Value = "  -4.65%  "
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  +5.74%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "6.32"
$ws.Range("E51").Value = "  +3.52%  "

# Restore original (default) style on column D now that the
# text values are locked in, so no stray formatting diff remains.
$dRange.Style = "Normal"
